# Update "want to go" counts (column F) across sheets to reflect the
# latest scrape snapshot, per commit "Update gh-pages to output generated
# at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 374
$ws1.Range("F5").Value = 1283
$ws1.Range("F6").Value = 215
$ws1.Range("F7").Value = 2456
$ws1.Range("F9").Value = 18482
$ws1.Range("F11").Value = 1865
$ws1.Range("F13").Value = 595
$ws1.Range("F14").Value = 318
$ws1.Range("F20").Value = 161

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 117

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5857

# Sheet "全部类型" (All types - aggregate of the above)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5857
$ws4.Range("F6").Value = 374
$ws4.Range("F10").Value = 1283
$ws4.Range("F12").Value = 215
$ws4.Range("F15").Value = 2456
$ws4.Range("F17").Value = 18482
$ws4.Range("F20").Value = 117
$ws4.Range("F21").Value = 117
$ws4.Range("F22").Value = 1865
$ws4.Range("F25").Value = 595
$ws4.Range("F26").Value = 318
$ws4.Range("F37").Value = 161
